# Add new measurement rows (2025-01-10 / serial 45667) to "adp", "bio_lucas"
# and "medidas" sheets, re-sort "medidas" chronologically, and refresh the
# sheet view selections, matching the upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "adp" (sheet1): append rows 98-113 for date 45667
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("adp")

$adp = @(
    @("dados", "peso", 79.75),
    @("dados", "triceps", 3),
    @("dados", "escapular", 10),
    @("dados", "biceps", 4),
    @("dados", "linha media axilar", 4),
    @("dados", "peitoral", 4),
    @("dados", "crista", 3),
    @("dados", "abdominal", 8),
    @("dados", "coxa", 5),
    @("dados", "panturrilha", 5),
    @("Jackson Pollock", "peso gordura", 7.29),
    @("Jackson Pollock", "peso magro", 72.46),
    @("Jackson Pollock", "não gordura", 90.86),
    @("Jackson Pollock", "gordura", 9.14),
    @("Jackson Pollock", "ideal gordura", 15.35),
    @("Jackson Pollock", "gordura sobra", -6.21)
)

$r = 98
foreach ($item in $adp) {
    $ws1.Cells.Item(2, 1).Copy($ws1.Cells.Item($r, 1))
    $ws1.Cells.Item($r, 1).Value = 45667
    $ws1.Cells.Item($r, 2).Value = $item[0]
    $ws1.Cells.Item($r, 3).Value = $item[1]
    $ws1.Cells.Item($r, 4).Value = $item[2]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "bio_lucas" (sheet2): append rows 77-91 for date 45667
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("bio_lucas")

$bio = @(
    @("dados do paciente", "massa adiposa", 10.13),
    @("dados do paciente", "Taxa metabolica", 1997),
    @("dados do paciente", "sarcopenix", 22),
    @("dados do paciente", "massa não adiposa", 69.62),
    @("dados do paciente", "massa muscular esqueletica", 39.72),
    @("dados do paciente", "Peso", 79.75),
    @("dados do paciente", "massa óssea", 3.4),
    @("dados do paciente", "idade metabolica", 28),
    @("dados do paciente", "Physique rating", 6),
    @("gráficos", "gordura corporal", 12.7),
    @("gráficos", "IMC", 26.04),
    @("gráficos", "Gordura visceral", 5.5),
    @("gráficos", "Analise Massa muscular", 83.01),
    @("gráficos", "Água corporal", 66.3),
    @("gráficos", "qualidade muscular", 68)
)

$r = 77
foreach ($item in $bio) {
    $ws2.Cells.Item(2, 1).Copy($ws2.Cells.Item($r, 1))
    $ws2.Cells.Item($r, 1).Value = 45667
    $ws2.Cells.Item($r, 2).Value = $item[0]
    $ws2.Cells.Item($r, 3).Value = $item[1]
    $ws2.Cells.Item($r, 4).Value = $item[2]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "medidas" (sheet3): sort existing rows chronologically by date, then
# append new rows 72-81 for date 45667
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("medidas")

$sortKey = $ws3.Range("A2:A71")
$ws3.Range("A2:C71").Sort($sortKey, 1)

$medidas = @(
    @("braço direito", 35.6),
    @("braço esquerdo", 35.6),
    @("antebraço direito", 35.6),
    @("antebraço esquerdo", 29.5),
    @("cintura", 77.5),
    @("quadril", 97.5),
    @("coxa esquerda", 58.5),
    @("coxa direita", 58.5),
    @("panturrilha direita", 38),
    @("panturrilha esquerda", 39.5)
)

$r = 72
foreach ($item in $medidas) {
    $ws3.Cells.Item(2, 1).Copy($ws3.Cells.Item($r, 1))
    $ws3.Cells.Item($r, 1).Value = 45667
    $ws3.Cells.Item($r, 2).Value = $item[0]
    $ws3.Cells.Item($r, 3).Value = $item[1]
    $r++
}

# Record the new sort state on the autoFilter (Excel does this automatically
# when you sort a filtered range from the UI).
$ws3.AutoFilter.Sort.SortFields.Clear()
$ws3.AutoFilter.Sort.SortFields.Add($ws3.Range("A2:A71"))
$ws3.AutoFilter.Sort.Apply()

# ---------------------------------------------------------------------------
# Restore/update the on-screen selections for each touched sheet, finishing
# on "inbody_full" so it remains the active tab (matches workbook.xml).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C116").Select()

$ws2.Activate()
$ws2.Range("D92").Select()

$ws3.Activate()
$ws3.Range("C82").Select()

$ws5 = $wb.Worksheets.Item("inbody_full")
$ws5.Activate()
$ws5.Range("E342").Select()
